$d = $word.ActiveDocument

# The whole (only) paragraph currently reads:
#   -Sí, y una polla- dijo Teo y se fue por la puerta.
# Replace it with "Zzzz", wrapped in proofErr spell-check markers, and
# leave a new trailing empty paragraph behind (mirrors Word inserting a
# paragraph mark after typing "Zzzz" + Enter and flagging it as a
# possible spelling mistake).

$r = $d.Paragraphs(1).Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:proofErr w:type="spellStart"/>' + `
              '<w:r><w:t>Zzzz</w:t></w:r>' + `
              '<w:proofErr w:type="spellEnd"/>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$r.InsertXML($xml)
